$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-03 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-04 Sunday", 2) | Out-Null
$d.Content.Find.Execute("94×24=2256", $true, $false, $false, $false, $false, $true, 1, $false, "33×55=1815", 2) | Out-Null
$d.Content.Find.Execute("98×92=9016", $true, $false, $false, $false, $false, $true, 1, $false, "77×68=5236", 2) | Out-Null
$d.Content.Find.Execute("91×86=7826", $true, $false, $false, $false, $false, $true, 1, $false, "46×28=1288", 2) | Out-Null
$d.Content.Find.Execute("13×63=819", $true, $false, $false, $false, $false, $true, 1, $false, "63×51=3213", 2) | Out-Null
$d.Content.Find.Execute("30×55=1650", $true, $false, $false, $false, $false, $true, 1, $false, "90×69=6210", 2) | Out-Null
$d.Content.Find.Execute("38×46=1748", $true, $false, $false, $false, $false, $true, 1, $false, "86×28=2408", 2) | Out-Null
$d.Content.Find.Execute("87×28=2436", $true, $false, $false, $false, $false, $true, 1, $false, "22×38=836", 2) | Out-Null
$d.Content.Find.Execute("86×59=5074", $true, $false, $false, $false, $false, $true, 1, $false, "28×87=2436", 2) | Out-Null
$d.Content.Find.Execute("51×33=1683", $true, $false, $false, $false, $false, $true, 1, $false, "87×36=3132", 2) | Out-Null
$d.Content.Find.Execute("99×91=9009", $true, $false, $false, $false, $false, $true, 1, $false, "67×62=4154", 2) | Out-Null
$d.Content.Find.Execute("52×73=3796", $true, $false, $false, $false, $false, $true, 1, $false, "56×88=4928", 2) | Out-Null
$d.Content.Find.Execute("27×45=1215", $true, $false, $false, $false, $false, $true, 1, $false, "79×52=4108", 2) | Out-Null
$d.Content.Find.Execute("91×91=8281", $true, $false, $false, $false, $false, $true, 1, $false, "98×69=6762", 2) | Out-Null
$d.Content.Find.Execute("49×41=2009", $true, $false, $false, $false, $false, $true, 1, $false, "15×91=1365", 2) | Out-Null
$d.Content.Find.Execute("58×64=3712", $true, $false, $false, $false, $false, $true, 1, $false, "54×82=4428", 2) | Out-Null
$d.Content.Find.Execute("96×99=9504", $true, $false, $false, $false, $false, $true, 1, $false, "55×13=715", 2) | Out-Null
$d.Content.Find.Execute("51×77=3927", $true, $false, $false, $false, $false, $true, 1, $false, "60×58=3480", 2) | Out-Null
$d.Content.Find.Execute("43×72=3096", $true, $false, $false, $false, $false, $true, 1, $false, "28×98=2744", 2) | Out-Null
$d.Content.Find.Execute("44×86=3784", $true, $false, $false, $false, $false, $true, 1, $false, "85×85=7225", 2) | Out-Null
$d.Content.Find.Execute("38×79=3002", $true, $false, $false, $false, $false, $true, 1, $false, "26×49=1274", 2) | Out-Null
$d.Content.Find.Execute("23×60=1380", $true, $false, $false, $false, $false, $true, 1, $false, "76×78=5928", 2) | Out-Null
$d.Content.Find.Execute("82×74=6068", $true, $false, $false, $false, $false, $true, 1, $false, "22×75=1650", 2) | Out-Null
$d.Content.Find.Execute("39×60=2340", $true, $false, $false, $false, $false, $true, 1, $false, "13×12=156", 2) | Out-Null
$d.Content.Find.Execute("76×23=1748", $true, $false, $false, $false, $false, $true, 1, $false, "33×23=759", 2) | Out-Null
$d.Content.Find.Execute("85×70=5950", $true, $false, $false, $false, $false, $true, 1, $false, "80×62=4960", 2) | Out-Null
